# Update "想去人数" (wanted-to-go count) figures across the sheets that
# track event attendance interest: 展览 (Exhibitions), 演出 (Performances),
# and 全部类型 (All types, an aggregate view of the two).
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 7654
$wsExpo.Range("F4").Value = 30
$wsExpo.Range("F5").Value = 466
$wsExpo.Range("F6").Value = 4318
$wsExpo.Range("F8").Value = 604
$wsExpo.Range("F11").Value = 162

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 12

# --- Sheet "全部类型" (All types, aggregate) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7654
$wsAll.Range("F5").Value = 30
$wsAll.Range("F6").Value = 466
$wsAll.Range("F7").Value = 4318
$wsAll.Range("F9").Value = 604
$wsAll.Range("F12").Value = 12
$wsAll.Range("F13").Value = 162

$wb.Save()
